$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.340.55"
$ws.Range("E2").Value = "  +3.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.970.88"
$ws.Range("E3").Value = "  +2.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.05"
$ws.Range("E5").Value = "  +1.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.73"
$ws.Range("E6").Value = "  +5.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.971.27"
$ws.Range("E8").Value = "  +2.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.00"
$ws.Range("E10").Value = "  +4.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  +2.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  +3.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.45"
$ws.Range("E14").Value = "  +6.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.126"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.461.87"
$ws.Range("E16").Value = "  +2.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.255.69"
$ws.Range("E17").Value = "  +3.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.89"
$ws.Range("E18").Value = "  +3.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.967.78"
$ws.Range("E19").Value = "  +1.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "445.04"
$ws.Range("E20").Value = "  +1.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.676"
$ws.Range("E22").Value = "  +2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.20"
$ws.Range("E23").Value = "  +4.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.50"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.96"
$ws.Range("E25").Value = "  +7.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.33"
$ws.Range("E26").Value = "  +2.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.19"
$ws.Range("E27").Value = "  +6.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  +8.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000109"
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.18"
$ws.Range("E31").Value = "  +4.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.56"
$ws.Range("E32").Value = "  +1.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.50"
$ws.Range("E34").Value = "  +3.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.975"
$ws.Range("E36").Value = "  +1.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  +3.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("E38").Value = "  +3.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.10"
$ws.Range("E39").Value = "  +6.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.01"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "43.74"
$ws.Range("E41").Value = "  +12.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.119"
$ws.Range("E42").Value = "  +1.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("E43").Value = "  +8.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.34"
$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "383.78"
$ws.Range("E45").Value = "  +12.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.761.35"
$ws.Range("E46").Value = "  +2.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0348"
$ws.Range("E47").Value = "  +3.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.47"
$ws.Range("E48").Value = "  +0.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000220"
$ws.Range("E50").Value = "  +9.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("E51").Value = "  +2.39%  "
